# forgot-password.xlsx: add a third "type" column that classifies each
# string row (title / description / input / button) next to the existing
# "string" (i18n key) and "value" (translated text) columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "type"
$ws.Range("C2").Value = "title"
$ws.Range("C3").Value = "description"
$ws.Range("C4").Value = "Email address input"
$ws.Range("C5").Value = "Button"
$ws.Range("C6").Value = "Button"

# Widen columns A and B to fit their (now taller) content, matching the
# author's resize of the key/value columns after the edit.
$ws.Columns.Item(1).ColumnWidth = 24.1666666666667
$ws.Columns.Item(2).ColumnWidth = 33.8307291666667

# Leave the selection where the author left it after entering the last
# value, one row below the new data.
$ws.Range("C7").Select()
